# Commit: feat: add 2022-Q1 data
#
# Insert a new "2022-Q1" worksheet (fund-holdings detail for 603712) right
# before the "总计" (summary) sheet, and prepend the new quarter's roll-up
# row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Reproduces the bold / centered / thin-bordered look used by every other
# quarter sheet's header row and index column (column A).
function Set-HeaderStyle($rng) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
    $rng.Borders.LineStyle = 1
}

# Writes a numeric-looking value (fund code, percentages, NAV, …) as TEXT —
# matching the source sheets, where these columns are inline strings, not
# numbers — by leading with an apostrophe text qualifier (preserves leading
# zeros in fund codes, e.g. "001475", and avoids silent numeric coercion).
function Set-TextValue($ws, $r, $c, $val) {
    $ws.Cells.Item($r, $c).Value = "'" + $val
}

# ---------------------------------------------------------------------------
# 1. Create the "2022-Q1" worksheet, positioned immediately before "总计".
#    NOTE: worksheet object refs resolve by position in this host, so
#    inserting a sheet shifts the index of every sheet after it. Re-fetch
#    "总计" by name AFTER the insert rather than reusing the pre-insert ref.
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$newSheet.Name = "2022-Q1"
$totalSheet = $wb.Worksheets.Item("总计")

# Header row (B1:H1)
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $cell = $newSheet.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 2]
    Set-HeaderStyle $cell
}

# Data rows (A2:H37): col A = row index (number, styled), B = fund code
# (text), C = fund name (text), D = fund size (text), E = total stock
# position (text), F = position share (text), G = held value ¥100M (text),
# H = position rank (number).
$data2022Q1 = @(
    @(0, '001475', '易方达国防军工混合', '180.12', '92.18', '6.03', '10.8612', 6),
    @(1, '010680', '华夏新兴成长股票A', '55.80', '87.37', '9.26', '5.1671', 1),
    @(2, '002251', '华夏军工安全灵活配置混合', '44.78', '94.71', '8.42', '3.7705', 3),
    @(3, '010305', '华夏创新驱动混合A', '28.56', '88.57', '8.81', '2.5161', 1),
    @(4, '012428', '华夏核心制造混合型证券投资基金A', '43.75', '91.55', '4.56', '1.9950', 8),
    @(5, '004698', '博时军工主题股票', '45.05', '89.94', '4.19', '1.8876', 10),
    @(6, '000001', '华夏成长混合', '31.69', '72.21', '4.23', '1.3405', 3),
    @(7, '002229', '华夏经济转型股票', '11.71', '86.14', '9.05', '1.0598', 1),
    @(8, '519908', '华夏兴华混合A', '9.39', '91.83', '8.91', '0.8366', 1),
    @(9, '960004', '华夏兴华混合H', '9.39', '91.83', '8.91', '0.8366', 1),
    @(10, '378010', '上投摩根成长先锋混合', '20.36', '83.15', '3.05', '0.6210', 6),
    @(11, '010681', '华夏新兴成长股票C', '5.29', '87.37', '9.26', '0.4899', 1),
    @(12, '006868', '华夏科技成长股票', '5.23', '87.33', '9.24', '0.4833', 1),
    @(13, '012429', '华夏核心制造混合型证券投资基金C', '9.10', '91.55', '4.56', '0.4150', 8),
    @(14, '013091', '上投摩根均衡优选混合A', '8.69', '64.17', '2.71', '0.2355', 5),
    @(15, '004640', '华夏节能环保股票', '5.45', '90.73', '3.72', '0.2027', 9),
    @(16, '010306', '华夏创新驱动混合C', '2.15', '88.57', '8.81', '0.1894', 1),
    @(17, '002703', '长城久源灵活配置混合', '1.56', '89.80', '7.41', '0.1156', 2),
    @(18, '000866', '华宝高端制造股票', '2.25', '91.51', '3.36', '0.0756', 6),
    @(19, '006952', '中银景元回报混合', '3.58', '33.73', '1.69', '0.0605', 3),
    @(20, '163823', '中银稳健策略灵活配置混合', '2.55', '48.41', '2.22', '0.0566', 6),
    @(21, '002067', '诺安精选回报灵活配置混合', '4.89', '25.39', '1.10', '0.0538', 5),
    @(22, '013899', '上投摩根全景优势股票A', '3.32', '46.28', '1.62', '0.0538', 9),
    @(23, '008773', '中银景泰回报混合', '4.83', '25.44', '1.07', '0.0517', 5),
    @(24, '002535', '中银鑫利灵活配置混合A', '6.75', '20.93', '0.75', '0.0506', 4),
    @(25, '003243', '上投摩根中国世纪灵活配置混合人民币份额（QDII）', '1.36', '84.74', '3.48', '0.0473', 4),
    @(26, '003244', '上投摩根中国世纪灵活配置混合美元现钞（QDII）', '1.36', '84.74', '3.48', '0.0473', 4),
    @(27, '003245', '上投摩根中国世纪灵活配置混合美元现汇（QDII）', '1.36', '84.74', '3.48', '0.0473', 4),
    @(28, '013092', '上投摩根均衡优选混合C', '0.77', '64.17', '2.71', '0.0209', 5),
    @(29, '002145', '诺安景鑫灵活配置混合', '0.53', '83.45', '3.65', '0.0193', 9),
    @(30, '002536', '中银鑫利灵活配置混合C', '2.30', '20.93', '0.75', '0.0172', 4),
    @(31, '002288', '中银稳进策略灵活配置混合', '0.70', '66.42', '2.37', '0.0166', 6),
    @(32, '004284', '华宝新优选一年定期开放灵活配置混合', '0.64', '38.91', '2.55', '0.0163', 1),
    @(33, '006890', '上投摩根领先优选混合', '0.36', '79.50', '3.25', '0.0117', 6),
    @(34, '320016', '诺安多策略混合', '0.19', '80.02', '3.33', '0.0063', 10),
    @(35, '013900', '上投摩根全景优势股票C', '0.20', '46.28', '1.62', '0.0032', 9),
)

for ($i = 0; $i -lt $data2022Q1.Length; $i++) {
    $row = $data2022Q1[$i]
    $r = $i + 2

    $cellA = $newSheet.Cells.Item($r, 1)
    $cellA.Value = $row[0]
    Set-HeaderStyle $cellA

    Set-TextValue $newSheet $r 2 $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    Set-TextValue $newSheet $r 4 $row[3]
    Set-TextValue $newSheet $r 5 $row[4]
    Set-TextValue $newSheet $r 6 $row[5]
    Set-TextValue $newSheet $r 7 $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------------
# 2. Rebuild the "总计" sheet, prepending the new 2022-Q1 summary row.
# ---------------------------------------------------------------------------
$totalSheet.Cells.Clear()

$totalSheet.Cells.Item(1, 2).Value = "日期"
$totalSheet.Cells.Item(1, 3).Value = "持有数量(只)"
$totalSheet.Cells.Item(1, 4).Value = "持有市值(亿元)"
Set-HeaderStyle $totalSheet.Range("B1:D1")

$dataTotal = @(
    @(0, '2022-Q1', 36, 33.68),
    @(1, '2021-Q4', 28, 38.6),
    @(2, '2021-Q3', 24, 32.45),
    @(3, '2021-Q2', 24, 30.01),
    @(4, '2021-Q1', 15, 16.25),
    @(5, '2020-Q4', 15, 18.26),
)

for ($i = 0; $i -lt $dataTotal.Length; $i++) {
    $row = $dataTotal[$i]
    $r = $i + 2

    $cellA = $totalSheet.Cells.Item($r, 1)
    $cellA.Value = $row[0]
    Set-HeaderStyle $cellA

    $totalSheet.Cells.Item($r, 2).Value = $row[1]
    $totalSheet.Cells.Item($r, 3).Value = $row[2]
    $totalSheet.Cells.Item($r, 4).Value = $row[3]
}

Write-Output "2022-Q1 sheet added; 总计 sheet updated."
